$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that currently sits at the end
#    of the "...password)" paragraph (right after the POST /login
#    section). It is a hidden/system bookmark so it never shows up in
#    Bookmarks.Count, but it is still addressable by name.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldBookmark = $d.Bookmarks("_GoBack")
    $oldBookmark.Delete()
}

# ------------------------------------------------------------------
# 2) Find the "/newLot ... creates a new lot" paragraph and append a
#    new italic run " (returns updated account)" to the end of it.
# ------------------------------------------------------------------
$target = $null
foreach ($para in $d.Paragraphs) {
    $t = $para.Range.Text
    if ($t -like "*/newLot*" -and $t -like "*creates a new lot*") {
        $target = $para
    }
}

$targetRange = $target.Range
# End of the paragraph's *text* (i.e. just before the paragraph mark).
$insertPos = $targetRange.End - 1

$insertPoint = $d.Range($insertPos, $insertPos)
$newText = " (returns updated account)"
$insertPoint.InsertAfter($newText)

$newRunRange = $d.Range($insertPos, $insertPos + $newText.Length)
$newRunRange.Font.Name = "Source Code Pro"
$newRunRange.Font.Italic = $true
$newRunRange.Font.Size = 11

# ------------------------------------------------------------------
# 3) Re-create the "_GoBack" bookmark, now collapsed at the very end
#    of that same paragraph's text (after the new run, before the
#    paragraph mark). Word's Bookmarks.Add needs a non-empty range to
#    respect the requested position, so temporarily wrap a throw-away
#    character, add the bookmark around it, then delete the
#    character -- deleting text strictly between a bookmark's start
#    and end collapses the bookmark in place instead of removing it.
# ------------------------------------------------------------------
$newEndPos = $insertPos + $newText.Length

$placeholderPoint = $d.Range($newEndPos, $newEndPos)
$placeholderPoint.InsertAfter("X")

$wrapRange = $d.Range($newEndPos, $newEndPos + 1)
$d.Bookmarks.Add("_GoBack", $wrapRange)

$placeholderRange = $d.Range($newEndPos, $newEndPos + 1)
$placeholderRange.Delete()
